$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.246.37'
$ws.Range("E2").Value = '  -2.30%  '
$ws.Range("D3").Value = '2.584.15'
$ws.Range("E3").Value = '  -2.56%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.77'
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.76'
$ws.Range("E6").Value = '  -3.09%  '
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("E8").Value = '  -1.68%  '
$ws.Range("D9").Value = '2.592.33'
$ws.Range("E9").Value = '  -3.30%  '
$ws.Range("E10").Value = '  -3.32%  '
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.159'
$ws.Range("E12").Value = '  +10.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("E13").Value = '  +3.20%  '
$ws.Range("D14").Value = '3.040.34'
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").Value = '59.208.19'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.96'
$ws.Range("E16").Value = '  +5.14%  '
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").Value = '2.588.61'
$ws.Range("E18").Value = '  -3.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.58'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '336.97'
$ws.Range("E20").Value = '  -2.23%  '
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.04'
$ws.Range("E24").Value = '  -3.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.466'
$ws.Range("E25").Value = '  +4.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("E27").Value = '  -2.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.36'
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  -3.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.93'
$ws.Range("E33").Value = '  +2.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.99'
$ws.Range("E34").Value = '  -1.64%  '
$ws.Range("E35").Value = '  -2.09%  '
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.878'
$ws.Range("E37").Value = '  -4.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.869'
$ws.Range("E38").Value = '  -5.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.47'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  -2.75%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '292.74'
$ws.Range("E42").Value = '  -5.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '132.52'
$ws.Range("E43").Value = '  +4.85%  '
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.953.73'
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.60'
$ws.Range("E51").Value = '  -1.84%  '
